$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7254747152328491
$ws.Range("B1").Value = 2.694529294967651
$ws.Range("C1").Value = 3.201138973236084
$ws.Range("D1").Value = 2.495908498764038
$ws.Range("E1").Value = 1.497864484786987
